# Refresh the cryptos list (prices / 1h volume, plus a few re-ranked rows)
# as published by the "Updated cryptos list ... with GitHub Actions" job.
#
# Price-column ("D") values are written with a leading apostrophe so Excel
# stores them as literal text (matching the sheet's existing inlineStr/text
# convention) instead of silently re-typing numeric-looking strings as
# numbers -- which would drop meaningful trailing zeros (e.g. "609.50" ->
# 609.5) or round very small values (e.g. "0.0000210" -> 0.000021).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''66.740.01'
$ws.Range('E2').Value = '  +0.93%  '
$ws.Range('D3').Value = '''3.622.87'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '''609.50'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').Value = '''149.90'
$ws.Range('E6').Value = '  +3.28%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '''0.490'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '''0.137'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').Value = '''8.03'
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('D11').Value = '''0.418'
$ws.Range('E11').Value = '  +1.34%  '
$ws.Range('D12').Value = '''4.217.98'
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').Value = '''0.0000210'
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('D14').Value = '''30.05'
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').Value = '''3.611.00'
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').Value = '''66.766.88'
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = '''11.69'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('D19').Value = '''6.39'
$ws.Range('E19').Value = '  +2.94%  '
$ws.Range('D20').Value = '''15.15'
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('D21').Value = '''432.77'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').Value = '''0.623'
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('D23').Value = '''78.82'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').Value = '''0.0000122'
$ws.Range('E25').Value = '  +2.73%  '
$ws.Range('D26').Value = '''8.40'
$ws.Range('E26').Value = '  +5.55%  '
$ws.Range('D27').Value = '''9.47'
$ws.Range('E27').Value = '  +3.56%  '
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = '''1.48'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('B31').Value = 'RenzoRestakedETH'
$ws.Range('C31').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D31').Value = '''3.605.58'
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '''25.66'
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = '''0.158'
$ws.Range('E33').Value = '  +3.50%  '
$ws.Range('D34').Value = '''7.94'
$ws.Range('E34').Value = '  +0.85%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '''5.70'
$ws.Range('E36').Value = '  +1.69%  '
$ws.Range('D37').Value = '''1.72'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('D38').Value = '''177.92'
$ws.Range('E38').Value = '  +2.05%  '
$ws.Range('D39').Value = '''0.0863'
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('D40').Value = '''5.28'
$ws.Range('E40').Value = '  +1.15%  '
$ws.Range('D41').Value = '''0.905'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').Value = '''1.91'
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '''45.88'
$ws.Range('E43').Value = '  -0.46%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '''2.58'
$ws.Range('E44').Value = '  +8.63%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '''0.998'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').Value = '''1.19'
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '''25.21'
$ws.Range('E47').Value = '  -2.54%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''24.08'
$ws.Range('E48').Value = '  +1.92%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = '''7.24'
$ws.Range('E49').Value = '  +1.52%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').Value = '''0.963'
$ws.Range('E50').Value = '  +1.98%  '
$ws.Range('D51').Value = '''0.238'
$ws.Range('E51').Value = '  -0.04%  '
